# Update "etat des virements" sheet: fix proprietaire (A:H) values and
# montants (I:J:K) for the existing 2 beneficiary rows, and append 3 new
# beneficiary rows (4-6) with their own account/amount data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns (A:H) to be stored as text, so numeric-looking
# values (account numbers, codes, etc.) are not coerced into numbers.
$ws.Range("A2:H6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "lala morale"
$ws.Range("B2").Value = "113564"
$ws.Range("C2").Value = "114321654687987654543213"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "23132"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/TEST DR/AV"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 24000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 24000

# Row 3
$ws.Range("A3").Value = "YASSINE TYEST"
$ws.Range("B3").Value = "BB125874"
$ws.Range("C3").Value = "115649679785432432321321"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "tesqt"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "001/TEST DR/AV"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 36000
$ws.Range("J3").Value = 5400
$ws.Range("K3").Value = 30600

# Row 4
$ws.Range("A4").Value = "Mustapha Tahiri"
$ws.Range("B4").Value = "BB147852"
$ws.Range("C4").Value = "114649797854321313211111"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "11111"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "002/tEST drrr/AV"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 20000
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 17000

# Row 5
$ws.Range("A5").Value = "Tawfiq mf"
$ws.Range("B5").Value = "BB169785"
$ws.Range("C5").Value = "114649778543212222222222"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "2222"
$ws.Range("F5").Value = "Logement de fonction"
$ws.Range("G5").Value = "002/LF/tEST drrr/AV"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 3000
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 2700

# Row 6
$ws.Range("A6").Value = "Test mf"
$ws.Range("B6").Value = "BB979797"
$ws.Range("C6").Value = "116649797855555555555555"
$ws.Range("D6").Value = "'"
$ws.Range("E6").Value = "1346"
$ws.Range("F6").Value = "Logement de fonction"
$ws.Range("G6").Value = "002/LF/tEST drrr/AV"
$ws.Range("H6").Value = "mensuelle"
$ws.Range("I6").Value = 7000
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 6300

